$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.572.68"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "3.509.13"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.05"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.32"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("D7").Value = "3.507.11"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("E10").Value = "  +3.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.64"
$ws.Range("E11").Value = "  +8.64%  "

$ws.Range("E12").Value = "  +1.84%  "

$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("D15").Value = "4.102.02"
$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").Value = "3.517.84"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").Value = "67.488.48"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("E19").Value = "  +2.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.63"
$ws.Range("E20").Value = "  +2.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.90"
$ws.Range("E21").Value = "  +6.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.17"
$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.17"
$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").Value = "3.648.76"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.77"
$ws.Range("E29").Value = "  +4.96%  "

$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.65"
$ws.Range("E31").Value = "  +7.05%  "

$ws.Range("E32").Value = "  +4.67%  "

$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.72"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("E35").Value = "  +0.90%  "

$ws.Range("E36").Value = "  +1.79%  "

$ws.Range("D37").Value = "3.504.16"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.04"
$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +7.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.48"
$ws.Range("E42").Value = "  -2.51%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0894"
$ws.Range("E43").Value = "  +2.56%  "

$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.41"
$ws.Range("E45").Value = "  +9.97%  "

$ws.Range("E46").Value = "  +0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.62"
$ws.Range("E47").Value = "  +2.61%  "

$ws.Range("E48").Value = "  +4.09%  "

$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("E50").Value = "  +1.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.19%  "
